# "PCA corrected + SFS start"
# Adds 16 new "PCA_model_dense_*" feature-importance columns (BM:CB) with
# headers in row 1, and a new data row (row 8) holding a new
# Sequential-Feature-Selection (SFS) run that populates the first few
# "normal" feature columns (A:J) plus the brand-new PCA columns (BM:CB).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New header cells BM1:CB1 - give them the same look (bold, centered,
#    bordered) as the rest of row 1 by copying A1's format first, then
#    stamp in the new header text (this also creates the 16 new shared
#    strings).
# ---------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("BM1:CB1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("BM1").Value = "PCA_model_dense_32_output_0 importance"
$ws.Range("BN1").Value = "PCA_model_dense_32_output_1 importance"
$ws.Range("BO1").Value = "PCA_model_dense_32_output_2 importance"
$ws.Range("BP1").Value = "PCA_model_dense_32_output_3 importance"
$ws.Range("BQ1").Value = "PCA_model_dense_32_output_4 importance"
$ws.Range("BR1").Value = "PCA_model_dense_32_output_5 importance"
$ws.Range("BS1").Value = "PCA_model_dense_32_output_6 importance"
$ws.Range("BT1").Value = "PCA_model_dense_32_output_7 importance"
$ws.Range("BU1").Value = "PCA_model_dense_32_output_8 importance"
$ws.Range("BV1").Value = "PCA_model_dense_32_output_9 importance"
$ws.Range("BW1").Value = "PCA_model_dense_32_output_10 importance"
$ws.Range("BX1").Value = "PCA_model_dense_32_output_11 importance"
$ws.Range("BY1").Value = "PCA_model_dense_16_output_0 importance"
$ws.Range("BZ1").Value = "PCA_model_dense_16_output_1 importance"
$ws.Range("CA1").Value = "PCA_model_dense_16_output_2 importance"
$ws.Range("CB1").Value = "PCA_model_dense_16_output_3 importance"

# ---------------------------------------------------------------------
# 2. New data row 8 - another study run.
# ---------------------------------------------------------------------
$ws.Range("A8").Value = 0.9999972175022197
$ws.Range("B8").Value = 10
$ws.Range("C8").Value = 0.5486301773758064
$ws.Range("D8").Value = 0.01029423128690497
$ws.Range("E8").Value = 0.5633201781999565
$ws.Range("F8").Value = 0.5298317981658773
$ws.Range("G8").Value = 0.06397796689902732
$ws.Range("H8").Value = 0.002266520493634809
$ws.Range("I8").Value = 0.03535079426427134
$ws.Range("J8").Value = 0.01050423920860266

$ws.Range("BM8").Value = 0.03984066441008055
$ws.Range("BN8").Value = 0.05516322053075325
$ws.Range("BO8").Value = 0.07962050867702015
$ws.Range("BP8").Value = 0.07461725168387211
$ws.Range("BQ8").Value = 0.06961474993311663
$ws.Range("BR8").Value = 0.03529056999609417
$ws.Range("BS8").Value = 0.03076375406550665
$ws.Range("BT8").Value = 0.04528849524931228
$ws.Range("BU8").Value = 0.03082261807805221
$ws.Range("BV8").Value = 0.03490381546225779
$ws.Range("BW8").Value = 0.04063722471022402
$ws.Range("BX8").Value = 0.04999402791620833
$ws.Range("BY8").Value = 0.06750081636007194
$ws.Range("BZ8").Value = 0.04942823645448381
$ws.Range("CA8").Value = 0.1131362200519838
$ws.Range("CB8").Value = 0.07127830555542623
